$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update province/ciudad names in column A (shared-string reordering in source diff
# translates to these cells displaying different text at the same row).
$ws.Range("A10").Value = 'Bizkaia/Vizcaya'
$ws.Range("A11").Value = 'Galicia'
$ws.Range("A13").Value = 'Valencia/Valencia'
$ws.Range("A14").Value = 'Aragon'
$ws.Range("A15").Value = 'Toledo'
$ws.Range("A16").Value = 'Navarra'
$ws.Range("A17").Value = 'Zaragoza'
$ws.Range("A19").Value = 'Araba/Alava'
$ws.Range("A20").Value = 'Valladolid'
$ws.Range("A21").Value = 'La Rioja'
$ws.Range("A22").Value = 'Salamanca'
$ws.Range("A23").Value = 'Alacant/Alicante'
$ws.Range("A24").Value = 'Malaga'
$ws.Range("A25").Value = 'Leon'
$ws.Range("A27").Value = 'Gipuzkoa/Guipuzcoa'
$ws.Range("A28").Value = 'Extremadura'
$ws.Range("A29").Value = 'Granada'
$ws.Range("A31").Value = 'Cuenca'
$ws.Range("A32").Value = 'Burgos'
$ws.Range("A33").Value = 'Asturias'
$ws.Range("A34").Value = 'Gran Canaria'
$ws.Range("A35").Value = 'Cantabria'
$ws.Range("A36").Value = 'Soria'
$ws.Range("A37").Value = 'Guadalajara'
$ws.Range("A38").Value = 'A Coruña'
$ws.Range("A39").Value = 'Caceres'
$ws.Range("A40").Value = 'Avila'
$ws.Range("A41").Value = 'Jaen'
$ws.Range("A42").Value = 'Pontevedra'
$ws.Range("A44").Value = 'Tenerife'
$ws.Range("A45").Value = 'Murcia'
$ws.Range("A46").Value = 'Castello/Castellon'
$ws.Range("A48").Value = 'Palencia'
$ws.Range("A49").Value = 'Huesca'
$ws.Range("A50").Value = 'Badajoz'
$ws.Range("A52").Value = 'Ourense'
$ws.Range("A53").Value = 'Almeria'
$ws.Range("A55").Value = 'Lugo'

# Update numeric data in columns B (Casos totales), C (Casos activos), D (Recuperados), E (Muertes)
$ws.Range("B4").Value = 64523
$ws.Range("C4").Value = 39227
$ws.Range("D4").Value = 16698
$ws.Range("E4").Value = 8598

$ws.Range("B5").Value = 51935
$ws.Range("C5").Value = 24063
$ws.Range("D5").Value = 22365
$ws.Range("E5").Value = 5507

$ws.Range("B6").Value = 17807
$ws.Range("C6").Value = 7317
$ws.Range("D6").Value = 8605
$ws.Range("E6").Value = 1885

$ws.Range("B7").Value = 16278
$ws.Range("C7").Value = 6062
$ws.Range("D7").Value = 7572
$ws.Range("E7").Value = 2738

$ws.Range("B8").Value = 13111
$ws.Range("C8").Value = 13985
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1403

$ws.Range("B9").Value = 12298
$ws.Range("C9").Value = 8779
$ws.Range("D9").Value = 2202
$ws.Range("E9").Value = 1317

$ws.Range("B10").Value = 9814
$ws.Range("C10").Value = 7124
$ws.Range("D10").Value = 4423
$ws.Range("E10").Value = 761

$ws.Range("B11").Value = 9216
$ws.Range("C11").Value = 7225
$ws.Range("D11").Value = 1400
$ws.Range("E11").Value = 591

$ws.Range("B12").Value = 8119
$ws.Range("E12").Value = 1022

$ws.Range("B13").Value = 5453
$ws.Range("C13").Value = 4066
$ws.Range("D13").Value = 2767
$ws.Range("E13").Value = 651

$ws.Range("B14").Value = 5291
$ws.Range("C14").Value = 3006
$ws.Range("D14").Value = 1462
$ws.Range("E14").Value = 823

$ws.Range("B15").Value = 5135
$ws.Range("C15").Value = 4178
$ws.Range("D15").Value = 10597
$ws.Range("E15").Value = 693

$ws.Range("B16").Value = 5039
$ws.Range("C16").Value = 2917
$ws.Range("D16").Value = 1635
$ws.Range("E16").Value = 487

$ws.Range("B17").Value = 4860
$ws.Range("C17").Value = 2261
$ws.Range("D17").Value = 1962
$ws.Range("E17").Value = 637

$ws.Range("B18").Value = 4799
$ws.Range("E18").Value = 473

$ws.Range("B19").Value = 4521
$ws.Range("C19").Value = 7124
$ws.Range("D19").Value = 4423
$ws.Range("E19").Value = 349

$ws.Range("B20").Value = 4145
$ws.Range("C20").Value = 1452
$ws.Range("D20").Value = 2351
$ws.Range("E20").Value = 342

$ws.Range("B21").Value = 4000
$ws.Range("C21").Value = 2581
$ws.Range("D21").Value = 1075
$ws.Range("E21").Value = 344

$ws.Range("B22").Value = 3792
$ws.Range("C22").Value = 1102
$ws.Range("D22").Value = 2346
$ws.Range("E22").Value = 344

$ws.Range("B23").Value = 3714
$ws.Range("C23").Value = 3017
$ws.Range("D23").Value = 1938
$ws.Range("E23").Value = 470

$ws.Range("B24").Value = 3703
$ws.Range("C24").Value = 1811
$ws.Range("D24").Value = 1626
$ws.Range("E24").Value = 266

$ws.Range("B25").Value = 3329
$ws.Range("C25").Value = 1502
$ws.Range("D25").Value = 1436
$ws.Range("E25").Value = 391

$ws.Range("B26").Value = 3193
$ws.Range("C26").Value = 828
$ws.Range("D26").Value = 2167
$ws.Range("E26").Value = 198

$ws.Range("B27").Value = 2951
$ws.Range("C27").Value = 7124
$ws.Range("D27").Value = 4423
$ws.Range("E27").Value = 273

$ws.Range("B28").Value = 2907
$ws.Range("C28").Value = 2354
$ws.Range("D28").Value = 76
$ws.Range("E28").Value = 477

$ws.Range("B29").Value = 2905
$ws.Range("C29").Value = 2126
$ws.Range("D29").Value = 508
$ws.Range("E29").Value = 271

$ws.Range("B30").Value = 2782
$ws.Range("C30").Value = 1341
$ws.Range("D30").Value = 1179
$ws.Range("E30").Value = 262

$ws.Range("B31").Value = 2733
$ws.Range("C31").Value = 4178
$ws.Range("D31").Value = 10597
$ws.Range("E31").Value = 291

$ws.Range("B32").Value = 2588
$ws.Range("C32").Value = 855
$ws.Range("D32").Value = 1530
$ws.Range("E32").Value = 203

$ws.Range("B33").Value = 2342
$ws.Range("C33").Value = 1012
$ws.Range("D33").Value = 1037
$ws.Range("E33").Value = 293

$ws.Range("B34").Value = 2250
$ws.Range("C34").Value = 1342
$ws.Range("D34").Value = 760
$ws.Range("E34").Value = 148

$ws.Range("B35").Value = 2241
$ws.Range("C35").Value = 1862
$ws.Range("D35").Value = 178
$ws.Range("E35").Value = 201

$ws.Range("B36").Value = 2149
$ws.Range("C36").Value = 369
$ws.Range("D36").Value = 1663
$ws.Range("E36").Value = 117

$ws.Range("B37").Value = 2122
$ws.Range("C37").Value = 4178
$ws.Range("D37").Value = 10597
$ws.Range("E37").Value = 234

$ws.Range("B38").Value = 1969
$ws.Range("C38").Value = 333
$ws.Range("D38").Value = 1788
$ws.Range("E38").Value = 67

$ws.Range("B39").Value = 1956
$ws.Range("C39").Value = 1316
$ws.Range("D39").Value = 260
$ws.Range("E39").Value = 380

$ws.Range("B40").Value = 1805
$ws.Range("C40").Value = 596
$ws.Range("D40").Value = 1078
$ws.Range("E40").Value = 131

$ws.Range("B41").Value = 1609
$ws.Range("C41").Value = 1001
$ws.Range("D41").Value = 444
$ws.Range("E41").Value = 164

$ws.Range("B42").Value = 1536
$ws.Range("C42").Value = 333
$ws.Range("D42").Value = 1411
$ws.Range("E42").Value = 30

$ws.Range("B43").Value = 1527
$ws.Range("C43").Value = 1093
$ws.Range("D43").Value = 330
$ws.Range("E43").Value = 104

$ws.Range("B44").Value = 1519
$ws.Range("C44").Value = 766
$ws.Range("D44").Value = 648
$ws.Range("E44").Value = 105

$ws.Range("B45").Value = 1506
$ws.Range("C45").Value = 1767
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 138

$ws.Range("B46").Value = 1447
$ws.Range("C46").Value = 1098
$ws.Range("D46").Value = 149
$ws.Range("E46").Value = 200

$ws.Range("B47").Value = 1406
$ws.Range("C47").Value = 412
$ws.Range("D47").Value = 856
$ws.Range("E47").Value = 138

$ws.Range("B48").Value = 1127
$ws.Range("C48").Value = 310
$ws.Range("D48").Value = 739
$ws.Range("E48").Value = 78

$ws.Range("B49").Value = 1009
$ws.Range("C49").Value = 320
$ws.Range("D49").Value = 592
$ws.Range("E49").Value = 97

$ws.Range("B50").Value = 944
$ws.Range("C50").Value = 997
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 92

$ws.Range("B51").Value = 798
$ws.Range("C51").Value = 303
$ws.Range("D51").Value = 414
$ws.Range("E51").Value = 81

$ws.Range("B52").Value = 751
$ws.Range("C52").Value = 333
$ws.Range("D52").Value = 660
$ws.Range("E52").Value = 22

$ws.Range("B53").Value = 618
$ws.Range("C53").Value = 376
$ws.Range("D53").Value = 193
$ws.Range("E53").Value = 49

$ws.Range("B54").Value = 611
$ws.Range("C54").Value = 336
$ws.Range("D54").Value = 195
$ws.Range("E54").Value = 80

$ws.Range("B55").Value = 586
$ws.Range("C55").Value = 333
$ws.Range("D55").Value = 520
$ws.Range("E55").Value = 11

$ws.Range("B56").Value = 500
$ws.Range("C56").Value = 275
$ws.Range("D56").Value = 180
$ws.Range("E56").Value = 45

$ws.Range("B60").Value = 78
$ws.Range("C60").Value = 38
$ws.Range("D60").Value = 37

$ws.Range("B61").Value = 71
$ws.Range("C61").Value = 56
$ws.Range("D61").Value = 11
$ws.Range("E61").Value = 4

$ws.Range("C63").Value = 21
$ws.Range("D63").Value = 3

$ws.Range("C67").Value = 7
$ws.Range("D67").Value = 0

# Update the timestamp text
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 13:04"
